$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains weekly price records for "Zapallo italiano" at the
# "Macroferia Regional de Talca" market. Two new weekly records are being
# added to the historical series:
#   - one inserted before the existing row 152 (original data shifts down)
#   - one inserted before the (shifted) original row 171 (remaining data
#     shifts down again)
# Net effect: the table grows from rows 2:182 to rows 2:184, with all
# existing rows preserved (just shifted) and two brand-new rows of data.

# --- Insert first new row at 152, pushing old 152:182 down to 153:183 ---
$ws.Rows.Item(152).Insert()

$ws.Range("A152").Value = 5
$ws.Range("B152").Value = "Macroferia Regional de Talca"
$ws.Range("C152").Value = "Maule"
$ws.Range("D152").Value = 44441
$ws.Range("E152").Value = 7
$ws.Range("F152").Value = 100112032
$ws.Range("G152").Value = "Zapallo italiano"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 300
$ws.Range("K152").Value = 15000
$ws.Range("L152").Value = 15000
$ws.Range("M152").Value = 15000
$ws.Range("N152").Value = "`$/caja 50 unidades"
$ws.Range("O152").Value = "Región de Arica y Parinacota"
$ws.Range("P152").Value = 300
$ws.Range("Q152").Value = 50
$ws.Range("R152").Value = "Hortaliza"

# --- Insert second new row at 172 (after old row 170, which is now at 171),
#     pushing the remaining rows (now 172:183) down to 173:184 ---
$ws.Rows.Item(172).Insert()

$ws.Range("A172").Value = 5
$ws.Range("B172").Value = "Macroferia Regional de Talca"
$ws.Range("C172").Value = "Maule"
$ws.Range("D172").Value = 44442
$ws.Range("E172").Value = 7
$ws.Range("F172").Value = 100112032
$ws.Range("G172").Value = "Zapallo italiano"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 250
$ws.Range("K172").Value = 15000
$ws.Range("L172").Value = 15000
$ws.Range("M172").Value = 15000
$ws.Range("N172").Value = "`$/caja 50 unidades"
$ws.Range("O172").Value = "Región de Arica y Parinacota"
$ws.Range("P172").Value = 300
$ws.Range("Q172").Value = 50
$ws.Range("R172").Value = "Hortaliza"
